# statistics.xlsx: shrink the "Statistics" table from a
# first_name/last_name/department/totalSum (4-col) sheet down to a
# first_name/email/gender (3-col) sheet, per "fix bag & remove logs".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The totalSum column (D) is dropped entirely.
$ws.Range("D1:D5").EntireColumn.Delete()

# Rows 3 and 4 (Nadiya/Samson data rows) are removed; only the header
# and the first data row remain, plus a leftover value cell at C3.
$ws.Range("A3:C4").ClearContents()

# New header row
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "gender"

# New data for row 2 (first_name stays "Alexia")
$ws.Range("B2").Value = "aalessandone18@clickbank.net"
$ws.Range("C2").Value = "Female"

# Leftover dollar-amount value, kept as literal text (not a number)
$ws.Range("C3").Value = "'$803.83"

# New column widths for the surviving columns (A is untouched)
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(3).ColumnWidth = 6
